# Add a new "qa_check3" column (column E) to the QA check report and
# update the fileid for the third data row (user input for capturing path).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the third QA check (copy D1's header formatting to E1)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "qa_check3"

# New column values (all "pass")
$ws.Range("E2").Value = "pass"
$ws.Range("E3").Value = "pass"
$ws.Range("E4").Value = "pass"

# Updated fileid for row 4 (captured path input changed 5 -> 6)
$ws.Range("A4").Value = 6
